# Swap the species-record data between row 13 and row 14, while leaving
# the shared/location columns (C, I, P, S, T, U, V, W, Y, AA, AD, AE, AG,
# AT, AW, AX, AY) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell13 = $ws.Range("$col`13")
    $cell14 = $ws.Range("$col`14")

    $v13 = $cell13.Value2
    $v14 = $cell14.Value2

    $cell13.Value = $v14
    $cell14.Value = $v13
}
